$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UQ")

# Update default glazing uncertainty values (column G, rows 11-19): 0.01 -> 0.04
for ($r = 11; $r -le 19; $r++) {
    $ws.Cells.Item($r, 7).Value = 0.04
}

# Update view state: active sheet selection / frozen pane scroll position
$ws.Activate()
$ws.Range("A3").Select()
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("E12").Select()
